# Atualização automática: 2025-08-08 22:00:24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: refined detection for the same fly (new image, tighter box, higher confidence)
$ws.Range("D16").Value = "image_20250807111314_ppp0.jpg"
$ws.Range("I16").Value = "'643,531,686,575"
$ws.Range("I16").Style = "Normal"
$ws.Range("J16").Value = "'0.76"
$ws.Range("J16").Style = "Normal"

# --- Row 17: refined detection for the same fly (new image, tighter box, higher confidence)
$ws.Range("D17").Value = "image_20250807111314_ppp0.jpg"
$ws.Range("I17").Value = "'794,481,830,526"
$ws.Range("I17").Style = "Normal"
$ws.Range("J17").Value = "'0.72"
$ws.Range("J17").Style = "Normal"

# --- Row 19: newly appended detection record
$ws.Range("A19").Value = "f77cad75-e373-4760-9d5a-1d927bfccd1d"
$ws.Range("B19").Value = "mosca"
$ws.Range("C19").Value = 45877
$ws.Range("C19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D19").Value = "image_20250808221835_ppp0.jpg"
$ws.Range("E19").Value = "PLACA_20250717165933"
$ws.Range("F19").Value = "Beja"
$ws.Range("G19").Value = 38.02035
$ws.Range("H19").Value = -7.94715
$ws.Range("I19").Value = "'819,160,858,215"
$ws.Range("I19").Style = "Normal"
$ws.Range("J19").Value = "'0.75"
$ws.Range("J19").Style = "Normal"
